# 0.9 Indexes, Variables/Constants (PL/SQL) & Constraints done.
# FUNCTIONS dictionary WIP.
#
# The "ADD / DROP CONSTRAINT" row is now finished, so its "Done" marker
# changes from the placeholder text "WIP" to the same numeric 1 used by
# every other completed row.
#
# The table also gets an AutoFilter on the "Done" column (C) that shows
# only the rows which are still blank (i.e. not yet done) - this hides
# every already-finished row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "ADD / DROP CONSTRAINT" (row 18) is now Done, same as the other
# finished rows -> numeric 1 instead of the "WIP" text.
$ws.Range("C18").Value = 1

# Filter column C (the 3rd column of A1:C18) to show blank cells only,
# via the classic checkbox-list form (Operator:=xlFilterValues,
# Criteria1:=Array("")), which hides every row whose "Done" cell is
# non-blank.
$ws.Range("A1:C18").AutoFilter(3, @(""), 7)

# The previous selection/scroll position (topLeftCell A7, cell C7) is no
# longer relevant once the sheet is filtered; the new selection is C8.
$ws.Range("C8").Select()
